$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 762
$ws.Range("I33").Value = 848.4167
$ws.Range("J33").Value = 243.5
$ws.Range("K33").Value = 848.4167
$ws.Range("L33").Value = 243.5
$ws.Range("M33").Value = -619.4167
$ws.Range("N33").Value = -701.5
$ws.Range("H38").Value = 501.26666
$ws.Range("I38").Value = 179.92857
$ws.Range("K38").Value = 539.78571
$ws.Range("M38").Value = -167.78571
$ws.Range("H53").Value = 395.2857
$ws.Range("I53").Value = 329.75
$ws.Range("J53").Value = 482.66666
$ws.Range("K53").Value = 329.75
$ws.Range("L53").Value = 482.66666
$ws.Range("M53").Value = 307.25
$ws.Range("N53").Value = -1756.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3973.5417
$ws.Range("I32").Value = 3711.5652
$ws.Range("K32").Value = 3711.5652
$ws.Range("M32").Value = -3424.5652
$ws.Range("H132").Value = 1123.6666
$ws.Range("I132").Value = 1123.6666
$ws.Range("K132").Value = 3370.9998
$ws.Range("M132").Value = -840.9998000000001
$ws.Range("H140").Value = 78214
$ws.Range("J140").Value = 78214
$ws.Range("L140").Value = 78214
$ws.Range("N140").Value = -88574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 354
$ws.Range("I64").Value = 98
$ws.Range("J64").Value = 439.33334
$ws.Range("K64").Value = 98
$ws.Range("L64").Value = 439.33334
$ws.Range("M64").Value = 127
$ws.Range("N64").Value = -889.33334
$ws.Range("H67").Value = 354
$ws.Range("I67").Value = 98
$ws.Range("J67").Value = 439.33334
$ws.Range("K67").Value = 98
$ws.Range("L67").Value = 439.33334
$ws.Range("M67").Value = 682
$ws.Range("N67").Value = -1999.33334
$ws.Range("H80").Value = 320.46667
$ws.Range("J80").Value = 371.16666
$ws.Range("L80").Value = 371.16666
$ws.Range("N80").Value = -2367.16666
$ws.Range("H83").Value = 320.46667
$ws.Range("J83").Value = 371.16666
$ws.Range("L83").Value = 1855.8333
$ws.Range("N83").Value = -11839.8333
$ws.Range("H94").Value = 1766.25
$ws.Range("J94").Value = 1500
$ws.Range("L94").Value = 1500
$ws.Range("N94").Value = -2402
$ws.Range("H105").Value = 4000
$ws.Range("I105").Value = 4000
$ws.Range("K105").Value = 4000
$ws.Range("M105").Value = -2253
$ws.Range("H107").Value = 5949.8335
$ws.Range("J107").Value = 7950
$ws.Range("L107").Value = 7950
$ws.Range("N107").Value = -11790
$ws.Range("H134").Value = 1599.5834
$ws.Range("I134").Value = 1567.7273
$ws.Range("K134").Value = 4703.1819
$ws.Range("M134").Value = -2168.1819
$ws.Range("H140").Value = 80640
$ws.Range("J140").Value = 80640
$ws.Range("L140").Value = 80640
$ws.Range("N140").Value = -91000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2080.8948
$ws.Range("I31").Value = 1823.8334
$ws.Range("K31").Value = 1823.8334
$ws.Range("M31").Value = -1528.8334
$ws.Range("H34").Value = 2080.8948
$ws.Range("I34").Value = 1823.8334
$ws.Range("K34").Value = 1823.8334
$ws.Range("M34").Value = -1621.8334
$ws.Range("H35").Value = 576.6667
$ws.Range("I35").Value = 576.6667
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 576.6667
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -282.6667
$ws.Range("N35").ClearContents()
$ws.Range("H58").Value = 1649.6666
$ws.Range("I58").Value = 1704.3334
$ws.Range("J58").Value = 1267
$ws.Range("K58").Value = 1704.3334
$ws.Range("L58").Value = 1267
$ws.Range("M58").Value = -1501.3334
$ws.Range("N58").Value = -1673
$ws.Range("H94").Value = 196505.83
$ws.Range("I94").Value = 232804.2
$ws.Range("K94").Value = 232804.2
$ws.Range("M94").Value = -232353.2
$ws.Range("H105").Value = 2830.2693
$ws.Range("J105").Value = 3500.4167
$ws.Range("L105").Value = 3500.4167
$ws.Range("N105").Value = -6994.4167
$ws.Range("H107").Value = 1244.5714
$ws.Range("I107").Value = 979.8
$ws.Range("K107").Value = 979.8
$ws.Range("M107").Value = 940.2
$ws.Range("H132").Value = 1649.862
$ws.Range("I132").Value = 1565.1923
$ws.Range("K132").Value = 4695.5769
$ws.Range("M132").Value = -2165.5769
$ws.Range("H134").Value = 971.6
$ws.Range("I134").Value = 989.3077
$ws.Range("K134").Value = 2967.9231
$ws.Range("M134").Value = -432.9231
$ws.Range("H136").Value = 1649.6666
$ws.Range("I136").Value = 1704.3334
$ws.Range("J136").Value = 1267
$ws.Range("K136").Value = 5113.0002
$ws.Range("L136").Value = 3801
$ws.Range("M136").Value = -2563.0002
$ws.Range("N136").Value = -8901

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 59953
$ws.Range("J37").Value = 59953
$ws.Range("L37").Value = 179859
$ws.Range("N37").Value = -180083
$ws.Range("H98").Value = 2339.8
$ws.Range("I98").Value = 198
$ws.Range("J98").Value = 2875.25
$ws.Range("K98").Value = 594
$ws.Range("L98").Value = 8625.75
$ws.Range("M98").Value = 904
$ws.Range("N98").Value = -11621.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2095.6667
$ws.Range("I122").Value = 2095.6667
$ws.Range("K122").Value = 6287.000100000001
$ws.Range("M122").Value = -3837.000100000001
$ws.Range("H132").Value = 1984.4375
$ws.Range("I132").Value = 1842.6154
$ws.Range("K132").Value = 5527.8462
$ws.Range("M132").Value = -2997.8462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1374.9286
$ws.Range("I22").Value = 1290
$ws.Range("J22").Value = 1587.25
$ws.Range("K22").Value = 1290
$ws.Range("L22").Value = 1587.25
$ws.Range("M22").Value = -995
$ws.Range("N22").Value = -2177.25
$ws.Range("H27").Value = 1374.9286
$ws.Range("I27").Value = 1290
$ws.Range("J27").Value = 1587.25
$ws.Range("K27").Value = 1290
$ws.Range("L27").Value = 1587.25
$ws.Range("M27").Value = -1183
$ws.Range("N27").Value = -1801.25
$ws.Range("H55").Value = 212.57143
$ws.Range("J55").Value = 298.66666
$ws.Range("L55").Value = 298.66666
$ws.Range("N55").Value = -644.66666
$ws.Range("H61").Value = 2587.25
$ws.Range("I61").Value = 2587.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2587.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2385.25
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 3000.0527
$ws.Range("I68").Value = 3000.1428
$ws.Range("J68").Value = 2999.8
$ws.Range("K68").Value = 3000.1428
$ws.Range("L68").Value = 2999.8
$ws.Range("M68").Value = -2251.1428
$ws.Range("N68").Value = -4497.8
$ws.Range("H71").Value = 3000.0527
$ws.Range("I71").Value = 3000.1428
$ws.Range("J71").Value = 2999.8
$ws.Range("K71").Value = 15000.714
$ws.Range("L71").Value = 14999
$ws.Range("M71").Value = -11256.714
$ws.Range("N71").Value = -22487
$ws.Range("H82").Value = 1811.5625
$ws.Range("I82").Value = 1868.2307
$ws.Range("J82").Value = 1566
$ws.Range("K82").Value = 1868.2307
$ws.Range("L82").Value = 1566
$ws.Range("M82").Value = -1507.2307
$ws.Range("N82").Value = -2288
$ws.Range("H85").Value = 1811.5625
$ws.Range("I85").Value = 1868.2307
$ws.Range("J85").Value = 1566
$ws.Range("K85").Value = 1868.2307
$ws.Range("L85").Value = 1566
$ws.Range("M85").Value = -620.2307000000001
$ws.Range("N85").Value = -4062
$ws.Range("H98").Value = 18000
$ws.Range("J98").Value = 18000
$ws.Range("L98").Value = 18000
$ws.Range("N98").Value = -23990
$ws.Range("H113").Value = 2587.25
$ws.Range("I113").Value = 2587.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2587.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -417.25
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 3596.5833
$ws.Range("I136").Value = 2715.5
$ws.Range("K136").Value = 8146.5
$ws.Range("M136").Value = -5596.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 30104
$ws.Range("J70").Value = 30104
$ws.Range("L70").Value = 30104
$ws.Range("N70").Value = -30734
$ws.Range("H73").Value = 30104
$ws.Range("J73").Value = 30104
$ws.Range("L73").Value = 30104
$ws.Range("N73").Value = -32288
$ws.Range("H132").Value = 3307
$ws.Range("I132").Value = 3724.2144
$ws.Range("J132").Value = 386.5
$ws.Range("K132").Value = 11172.6432
$ws.Range("L132").Value = 1159.5
$ws.Range("M132").Value = -8642.643199999999
$ws.Range("N132").Value = -6219.5
$ws.Range("H136").Value = 2508.6128
$ws.Range("I136").Value = 2440.276
$ws.Range("J136").Value = 3499.5
$ws.Range("K136").Value = 7320.828
$ws.Range("L136").Value = 10498.5
$ws.Range("M136").Value = -4770.828
$ws.Range("N136").Value = -15598.5
